$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.092.75'
$ws.Range('E2').Value = '  -1.85%  '

# Row 3
$ws.Range('D3').Value = '1.830.59'
$ws.Range('E3').Value = '  -3.24%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  -0.09%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '229.41'
$ws.Range('E5').Value = '  -4.01%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').Value = '  -0.06%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4630'
$ws.Range('E7').Value = '  -4.07%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2693'
$ws.Range('E8').Value = '  -6.25%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06212'
$ws.Range('E9').Value = '  -5.26%  '

# Row 10
$ws.Range('D10').Value = '1.826.89'
$ws.Range('E10').Value = '  -4.72%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07343'
$ws.Range('E11').Value = '  -1.73%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '16.00'
$ws.Range('E12').Value = '  -4.26%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.886'
$ws.Range('E13').Value = '  -4.35%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '82.84'
$ws.Range('E14').Value = '  -6.10%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6173'
$ws.Range('E15').Value = '  -7.69%  '

# Row 16
$ws.Range('D16').Value = '30.050.24'
$ws.Range('E16').Value = '  -1.93%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.9994'
$ws.Range('E17').Value = '  -0.13%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '226.13'
$ws.Range('E18').Value = '  -2.99%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007241'
$ws.Range('E19').Value = '  -4.48%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.05%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.26'
$ws.Range('E21').Value = '  -7.68%  '

# Row 22
$ws.Range('D22').Value = '2.068.32'
$ws.Range('E22').Value = '  -5.81%  '

# Row 23
$ws.Range('E23').Value = '  -8.95%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.809'
$ws.Range('E24').Value = '  -6.81%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '165.54'
$ws.Range('E25').Value = '  -2.20%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.105'
$ws.Range('E26').Value = '  -2.66%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '17.72'
$ws.Range('E27').Value = '  -6.08%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.838'
$ws.Range('E28').Value = '  -6.42%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.1016'
$ws.Range('E29').Value = '  -0.70%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.366'
$ws.Range('E30').Value = '  -2.10%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.035'
$ws.Range('E31').Value = '  -6.83%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.743'
$ws.Range('E32').Value = '  -7.26%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04789'
$ws.Range('E33').Value = '  -5.46%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.124'
$ws.Range('E34').Value = '  -7.57%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.6954'
$ws.Range('E35').Value = '  -7.90%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.679'
$ws.Range('E36').Value = '  -1.19%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.603'
$ws.Range('E38').Value = '  -1.68%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.8892'
$ws.Range('E39').Value = '  -3.30%  '

# Row 40
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9993'
$ws.Range('E40').Value = '  -0.36%  '

# Row 41
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.906'
$ws.Range('E41').Value = '  -8.12%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '103.10'
$ws.Range('E42').Value = '  -3.95%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.452'
$ws.Range('E43').Value = '  -3.93%  '

# Row 44
$ws.Range('E44').Value = '  -7.37%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '6.868'
$ws.Range('E45').Value = '  -7.81%  '

# Row 46
$ws.Range('E46').Value = '  -7.19%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '58.92'
$ws.Range('E47').Value = '  -8.87%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.452'
$ws.Range('E48').Value = '  -6.11%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.05519'
$ws.Range('E49').Value = '  -2.52%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '32.43'
$ws.Range('E50').Value = '  -4.42%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.348'
$ws.Range('E51').Value = '  -10.07%  '
